# Update the "Monthly Consumptions" summary report figures.
# Table 2 = "I. Total Cost per Section" summary table
# Table 3 = weekly ("WEEK 1".."WEEK 5"/"Total") breakdown table
#
# For both tables the Production row's Total (WEEK 4 + grand Total),
# the Meter Maintenance row's Total (WEEK 4 + grand Total), and the
# overall Total/TOTAL row's Total (WEEK 4 + grand Total) increase in
# magnitude.

$d = $word.ActiveDocument

function Set-CellText($cell, $newText) {
    $rng = $cell.Range
    # Cell.Range includes the trailing end-of-cell mark; trim it so we
    # only overwrite the cell's actual content and keep its formatting.
    $contentRng = $d.Range($rng.Start, $rng.End - 1)
    $contentRng.Text = $newText
}

# --- Table 2: "I. Total Cost per Section" ---
$tblSummary = $d.Tables.Item(2)
Set-CellText $tblSummary.Cell(6, 3) "-25240"    # b. Production
Set-CellText $tblSummary.Cell(7, 3) "-39830"    # c. Meter Maintenance
Set-CellText $tblSummary.Cell(13, 3) "-73740"   # Total

# --- Table 3: weekly breakdown ---
$tblWeekly = $d.Tables.Item(3)
Set-CellText $tblWeekly.Cell(3, 5) "-25240"     # Production, WEEK 4
Set-CellText $tblWeekly.Cell(3, 7) "-25240"     # Production, Total
Set-CellText $tblWeekly.Cell(4, 5) "-39830"     # Meter Maintenance, WEEK 4
Set-CellText $tblWeekly.Cell(4, 7) "-39830"     # Meter Maintenance, Total
Set-CellText $tblWeekly.Cell(10, 5) "-73740"    # TOTAL, WEEK 4
Set-CellText $tblWeekly.Cell(10, 7) "-73740"    # TOTAL, Total
